# pro_reports_page.pptx - "Changed design on reporting page"
#
# The annotation callout pair (yellow "Disaster ID" oval marker + its
# label textbox) on slide 1 shifts 9525 EMU (0.75pt / 1px @ 96dpi) to the
# left, matching a fine nudge of the pair in the PowerPoint UI.
#
# Target OOXML (EMU):
#   Oval 11      : x 8430936 -> 8421411   (y stays 868509)
#   TextBox 12   : x 8443520 -> 8433995   (y stays 872629)
#
# Shape.Left/Top are exposed by PowerPoint as Single (32-bit float)
# points, so the literals below are chosen to round-trip to the exact
# target EMU through that Single-precision points representation
# (point = EMU / 12700).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$oval = $s.Shapes.Item("Oval 11")
$oval.Left = 663.103271484375

$textLabel = $s.Shapes.Item("TextBox 12")
$textLabel.Left = 664.0941162109375
